# Update the master output for the "three-digit number divided by
# one-digit number" worksheet to the problems generated for 2025-07-24.

$d = $word.ActiveDocument

function Set-RangeText($range, [string]$newText) {
    # Replace the visible text of $range while leaving any trailing
    # paragraph mark / end-of-cell marker untouched.
    $scoped = $d.Range($range.Start, $range.End - 1)
    $scoped.Text = $newText
}

# --- Header date -----------------------------------------------------
Set-RangeText $d.Paragraphs(1).Range "2025-07-24 Thursday"

# --- Table of division problems --------------------------------------
$t = $d.Tables(1)

function Set-CellText($row, $col, [string]$newText) {
    Set-RangeText $t.Cell($row, $col).Range $newText
}

# Row 1
Set-CellText 1 1 "663÷4="
Set-CellText 1 2 "896÷8="
Set-CellText 1 3 "716÷2="
Set-CellText 1 4 "333÷5="
Set-CellText 1 5 "665÷6="

# Row 5
Set-CellText 5 1 "344÷7="
Set-CellText 5 2 "439÷6="
Set-CellText 5 3 "284÷8="
Set-CellText 5 4 "161÷3="
Set-CellText 5 5 "641÷8="

# Row 9
Set-CellText 9 1 "518÷2="
Set-CellText 9 2 "400÷2="
Set-CellText 9 3 "582÷5="
Set-CellText 9 4 "704÷6="
Set-CellText 9 5 "722÷3="

# Row 13
Set-CellText 13 1 "687÷6="
Set-CellText 13 2 "557÷9="
Set-CellText 13 3 "495÷3="
Set-CellText 13 4 "507÷3="
Set-CellText 13 5 "163÷8="

# Row 17
Set-CellText 17 1 "297÷2="
Set-CellText 17 2 "409÷6="
Set-CellText 17 3 "626÷6="
Set-CellText 17 4 "261÷2="
Set-CellText 17 5 "750÷3="

Write-Output "done"
